$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New CFP rows 162-167 appended to the "listCFPs" sheet.
# Shared strings are introduced in a specific order (mirrors the order they
# were first typed into the workbook), so cell values are assigned in that
# exact sequence to reproduce the resulting sharedStrings.xml table.
# ---------------------------------------------------------------------------

# --- Row 162 : mydefichain Ocean 2023 Cycle-2 ---
$ws.Range("B162").Value = "06/2023"
$ws.Range("C162").Value = "mydefichain Ocean 2023 Cycle-2"
$ws.Hyperlinks.Add($ws.Range("E162"), "https://www.reddit.com/r/defiblockchain/comments/12qhhj5/mydefichain_ocean_2023_actual_cycle_round/") | Out-Null
$ws.Range("E162").Style = "Link"
$ws.Range("D162").Value = "Bernd Mack and Andreas Lentz with mydefichain"
$ws.Range("A162").Value = 162
$ws.Range("F162").Value = 7500
$ws.Range("G162").Value = "passed"

# --- Row 163 : Telegram Moderators H1 2023 ---
$ws.Range("C163").Value = "Telegram Moderators H1 2023"
$ws.Hyperlinks.Add($ws.Range("E163"), "https://www.reddit.com/r/defiblockchain/comments/12rugem/telegram_moderators_h1_2023/") | Out-Null
$ws.Range("E163").Style = "Link"
$ws.Range("B163").Value = "06/2023"
$ws.Range("D163").Value = "Telegram moderators"
$ws.Range("A163").Value = 163
$ws.Range("F163").Value = 82750
$ws.Range("G163").Value = "passed"

# --- Row 164 : Integration of Yield Machine and Interest Rates ---
$ws.Range("C164").Value = "Integration of Yield Machine and Interest Rates on defichain-income.com"
$ws.Hyperlinks.Add($ws.Range("E164"), "https://www.reddit.com/r/defiblockchain/comments/125s2s7/integration_of_yield_machine_and_interest_rates/") | Out-Null
$ws.Range("E164").Style = "Link"
$ws.Range("B164").Value = "06/2023"
$ws.Range("D164").Value = "Igor Shelkovenkov"
$ws.Range("A164").Value = 164
$ws.Range("F164").Value = 10000
$ws.Range("G164").Value = "passed"

# --- Row 165 : mydefichain-Ocean-2023-cycle-3 ---
$ws.Range("C165").Value = "mydefichain-Ocean-2023-cycle-3"
$ws.Range("B165").Value = "07/2023"
$ws.Hyperlinks.Add($ws.Range("E165"), "https://www.reddit.com/r/defiblockchain/comments/12qhhj5/mydefichain_ocean_2023_actual_cycle_round/") | Out-Null
$ws.Range("E165").Style = "Link"
$ws.Range("D165").Value = "Bernd Mack and Andreas Lentz with mydefichain"
$ws.Range("A165").Value = 165
$ws.Range("F165").Value = 11000
$ws.Range("G165").Value = "passed"

# --- Row 166 : mydefichain-Ocean-2023-cycle-4 ---
$ws.Range("C166").Value = "mydefichain-Ocean-2023-cycle-4"
$ws.Range("B166").Value = "09/2023"
$ws.Hyperlinks.Add($ws.Range("E166"), "https://www.reddit.com/r/defiblockchain/comments/12qhhj5/mydefichain_ocean_2023_actual_cycle_round/") | Out-Null
$ws.Range("E166").Style = "Link"
$ws.Range("D166").Value = "Bernd Mack and Andreas Lentz with mydefichain"
$ws.Range("A166").Value = 166
$ws.Range("F166").Value = 11000
$ws.Range("G166").Value = "passed"

# --- Row 167 : DeFiChain Turkey ---
$ws.Range("C167").Value = "DeFiChain Turkey - On the Way to Becoming the Biggest Regional Community"
$ws.Hyperlinks.Add($ws.Range("E167"), "https://www.reddit.com/r/defiblockchain/comments/1456fnx/cfp_defichain_turkey_on_the_way_to_becoming_the/") | Out-Null
$ws.Range("E167").Style = "Link"
$ws.Range("D167").Value = "DeFiChain Turkey"
$ws.Range("B167").Value = "09/2023"
$ws.Range("A167").Value = 167
$ws.Range("F167").Value = 45000
$ws.Range("G167").Value = "passed"

# Update the current selection / active cell to reflect where the author
# ended up after entering the new data.
$ws.Range("E170").Select() | Out-Null
